$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for the "16-jun" column
$ws.Range("I1").Value = "16-jun"

# Fill in the new column I values (row 2 through row 11)
$values = @(15, 13, 6, 13, 15, 14, 12, 15, 17, 8)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

# Match the number formatting/alignment of the adjacent column H for the data cells
$ws.Range("H2:H11").Copy()
$ws.Range("I2:I11").PasteSpecial(-4122)

# Update selection to mirror the post-edit state
$ws.Range("I12").Select()
